# Update F-column (想去人数 / "want to go" counts) values across all sheets
$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 723
$ws.Range("F5").Value = 415
$ws.Range("F6").Value = 687
$ws.Range("F10").Value = 930
$ws.Range("F14").Value = 51
$ws.Range("F17").Value = 23983
$ws.Range("F18").Value = 2217
$ws.Range("F19").Value = 142
$ws.Range("F20").Value = 353
$ws.Range("F21").Value = 32
$ws.Range("F22").Value = 52
$ws.Range("F25").Value = 66
$ws.Range("F26").Value = 228
$ws.Range("F28").Value = 54
$ws.Range("F30").Value = 346
$ws.Range("F32").Value = 433

# Sheet: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F7").Value = 254
$ws.Range("F8").Value = 22
$ws.Range("F9").Value = 245
$ws.Range("F10").Value = 3599
$ws.Range("F13").Value = 3
$ws.Range("F15").Value = 29
$ws.Range("F16").Value = 18
$ws.Range("F19").Value = 4117

# Sheet: 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 162
$ws.Range("F4").Value = 749
$ws.Range("F5").Value = 240

# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 162
$ws.Range("F5").Value = 749
$ws.Range("F6").Value = 723
$ws.Range("F8").Value = 415
$ws.Range("F9").Value = 687
$ws.Range("F14").Value = 254
$ws.Range("F15").Value = 240
$ws.Range("F18").Value = 930
$ws.Range("F21").Value = 51
$ws.Range("F24").Value = 23983
$ws.Range("F25").Value = 22
$ws.Range("F26").Value = 245
$ws.Range("F29").Value = 3
$ws.Range("F30").Value = 2217
$ws.Range("F31").Value = 142
$ws.Range("F32").Value = 353
$ws.Range("F33").Value = 32
$ws.Range("F35").Value = 29
$ws.Range("F38").Value = 228
$ws.Range("F40").Value = 18
$ws.Range("F41").Value = 54
$ws.Range("F46").Value = 433
$ws.Range("F48").Value = 4117

